$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor data window: the series advanced by one sample (old row 2 dropped)
# and 10 new trailing samples were appended (rows 22-31), per the "May 9th" update.
$data = @(
    @(-0.4181592197061701, 4.999134939407615, 4.467728652686705),
    @(-0.2939812507584775, 1.309815360006858, 2.089699776373152),
    @(0.6610930461749918, -2.38056925078418, -1.609635165918661),
    @(1.810891748588788, -1.378563615763306, -2.978032114349805),
    @(2.31189947707626, -2.220112597830942, -3.509698769756517),
    @(-4.20486067611483, -2.35995082988918, 3.951543032566025),
    @(-6.209334919385752, -1.759963089060472, 10.33535447298917),
    @(1.966778721764825, 7.234104103016977, 2.021741985160567),
    @(2.757100457342962, -3.825466089159538, -3.918381726630348),
    @(6.505795770716528, -1.982219981255937, -9.447731205236206),
    @(1.019708882982564, -2.124502529607678, -4.567799980673901),
    @(-3.90443126063479, -6.289577069683729, 0.08947662911672927),
    @(-4.737582387211194, 1.679735990328154, 5.945494293052519),
    @(-2.958322148456745, 5.940943116339433, 1.412405640165503),
    @(-3.080484191948098, 14.34931452920528, -6.689891151178712),
    @(3.206049725274078, -2.906924434911418, -3.371373475154916),
    @(1.706316604792587, -3.271737626170626, -8.0454896521347),
    @(-3.738767730855505, 0.06012086237769707, 3.632486434740414),
    @(-3.090678747569279, 1.963554282054647, 5.732710114149367),
    @(-3.92490029335021, 8.097008705139087, -0.6271078586577614),
    @(-1.013831214370013, 11.75312601071651, -9.512279387946442),
    @(1.26375397343502, -1.259433144720926, -2.911592367653529),
    @(3.293022926722733, -6.641965772504008, -9.101227296847053),
    @(-0.6729963912028092, -6.046104440065187, 2.590652748803134),
    @(-4.304234186622589, 1.124698318053637, 5.37724533482133),
    @(-3.525356531143185, 1.932811015120178, -0.6847266883493761),
    @(0.02450294917989604, 11.26465672644501, -6.819475414596996),
    @(1.655811092563884, -3.081682441390621, -0.3115755630151993),
    @(4.203615091671459, -2.449124443196801, -8.744278280116706),
    @(-2.318640726749061, -4.338136873512635, 1.472857043007876)
)

$row = 2
foreach ($values in $data) {
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $row = $row + 1
}

